# service-cloning does now work
# Adds a "CLONEFROM" column (M) to the host table, sets testhost03's
# CLONEFROM value to "Linuxtest01", and clears the stray max_check_attempts
# value that had been left in K5 (testhost04 row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column M
$ws.Range("M1").Value = "CLONEFROM"

# testhost03 (row 4) now clones from Linuxtest01
$ws.Range("M4").Value = "Linuxtest01"

# Remove the leftover max_check_attempts value on testhost04 (row 5)
$ws.Range("K5").ClearContents()

# Match the author's final selection in the sheet view
$ws.Range("M5").Select()
